# Bosnia Herzegovina Premier Liga - base update (19-04-2024 00:38)
#
# 1) Four pairs of rows had their match data (everything except the
#    leading id label in column A) swapped between the two rows.
# 2) One new match (row 158) was appended at the end of the sheet.

function Set-RowValues {
    param($ws, $row, $startCol, $values)
    for ($i = 0; $i -lt $values.Count; $i++) {
        $v = $values[$i]
        if ($null -ne $v) {
            $ws.Cells.Item($row, $startCol + $i).Value = $v
        }
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap data between row 36 and row 37 ---
Set-RowValues $ws 36 2 @(6864629, "Bosnia Herzegovina Premier Liga", "Bosnia  Herzegovina Premier Liga", 45186.61458333334, "Borac Banja Luka", "NK Posusje", 1, 0, "H", 1.363, 4.5, 6.5, 1.363, 4.2, 6.5, -1.25, 1.95, 1.85, 2.5, 1.925, 1.875, 0.363, -1, -1, -0.5, 0.425, -1, 0.875)
Set-RowValues $ws 37 2 @(6865299, "Bosnia Herzegovina Premier Liga", "Bosnia  Herzegovina Premier Liga", 45186.61458333334, "Siroki Brijeg", "Zvijezda 09", 2, 1, "H", 1.25, 5.5, 8, 1.4, 4.75, 5.75, -1.25, 1.9, 1.9, 2.75, 1.85, 1.95, 0.3999999999999999, -1, -1, -0.5, 0.45, 0.425, -0.5)

# --- Swap data between row 49 and row 50 ---
Set-RowValues $ws 49 2 @(6865310, "Bosnia Herzegovina Premier Liga", "Bosnia  Herzegovina Premier Liga", 45200.41666666666, "NK Igman Konjic", "Zrinjski Mostar", 0, 2, "A", 3.4, 3.6, 1.833, 4.75, 4.75, 1.45, 1.25, 1.775, 2.025, 2.75, 1.85, 1.95, -1, -1, 0.45, -1, 1.025, -1, 0.95)
Set-RowValues $ws 50 2 @(6865311, "Bosnia Herzegovina Premier Liga", "Bosnia  Herzegovina Premier Liga", 45200.41666666666, "Sloga", "GOSK Gabela", 3, 2, "H", 1.833, 3.6, 3.4, 1.909, 3.4, 3.3, -0.5, 1.925, 1.875, 2.25, 1.825, 1.975, 0.909, -1, -1, 0.925, -1, 0.825, -1)

# --- Swap data between row 87 and row 88 ---
Set-RowValues $ws 87 2 @(7505495, "Bosnia Herzegovina Premier Liga", "Bosnia  Herzegovina Premier Liga", 45256.375, "Sloga", "Zvijezda 09", 1, 0, "H", 1.444, 4.2, 5.5, 1.5, 4.2, 5.25, -1, 1.8, 2, 2.75, 1.775, 2.025, 0.5, -1, -1, 0, 0, -1, 1.025)
Set-RowValues $ws 88 2 @(7505497, "Bosnia Herzegovina Premier Liga", "Bosnia  Herzegovina Premier Liga", 45256.375, "Zeljeznicar", "NK Posusje", 1, 1, "D", 1.65, 3.4, 4.75, 1.8, 3.2, 4.2, -0.5, 1.825, 1.975, 2, 1.75, 2.05, -1, 2.2, -1, -1, 0.9750000000000001, 0, 0)

# --- Swap data between row 99 and row 100 ---
Set-RowValues $ws 99 2 @(6865343, "Bosnia Herzegovina Premier Liga", "Bosnia  Herzegovina Premier Liga", 45269.375, "Sloga", "NK Posusje", 1, 0, "H", 1.909, 3.3, 3.5, 2.2, 2.8, 3.3, -0.25, 1.95, 1.85, 1.75, 1.875, 1.925, 1.2, -1, -1, 0.95, -1, -1, 0.925)
Set-RowValues $ws 100 2 @(6864639, "Bosnia Herzegovina Premier Liga", "Bosnia  Herzegovina Premier Liga", 45269.375, "Zvijezda 09", "Borac Banja Luka", 1, 2, "A", 11, 6, 1.2, 10, 6.5, 1.181, 2, 1.825, 1.975, 3, 1.9, 1.9, -1, -1, 0.181, 0.825, -1, 0, 0)

# --- Append new row 158 (upcoming fixture, no result yet) ---
# Copy formats from the last existing data row (157) for the styled cells
# (id column has a bold/bordered/centered style, date column has a custom
# date number format), then overwrite values.
$ws.Cells.Item(157, 1).Copy()
$ws.Cells.Item(158, 1).PasteSpecial(-4122)
$ws.Cells.Item(158, 1).Value = 156

$ws.Cells.Item(157, 5).Copy()
$ws.Cells.Item(158, 5).PasteSpecial(-4122)
$ws.Cells.Item(158, 5).Value = 45401.625

Set-RowValues $ws 158 2 @(7952750, "Bosnia Herzegovina Premier Liga", "Bosnia  Herzegovina Premier Liga", $null, "Velez Mostar", "Zvijezda 09", $null, $null, $null, 1.25, 5.75, 7, 1.2, 5.25, 12, -1.75, 1.825, 1.975, 3, 1.95, 1.85, 0, 0, 0, 0, 0, $null, $null)
